$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "Perioadele campaniei din Constelația Pegasus 2022:",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Perioadele campaniei din 2022 pentru Constelația Pegasus:",
    2
)
